# Updated symbol list on Sat Jan  7 11:34:58 UTC 2023 with GitHub Actions
# Applies updated Price (column D) and Volume(1h) (column E) values
# for the cryptos.xlsx sheet. All target cells are stored as text,
# so each cell's number format is forced to Text ("@") before the
# numeric-looking / percent-looking string is assigned, which keeps
# Excel from auto-converting the literal string into a number cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "260.90"

$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "1.88%"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "27.39"

$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "2.41%"

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "4.689"

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "0.12%"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.06091"

$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "2.54%"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "6.665"

$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "0.79%"

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.8461"

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "-0.66%"

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "1.03%"

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1406"

$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "1.95%"

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.04888"

$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "11.98%"

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07097"

$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "1.31%"

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.03076"

$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "0.35%"

$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-0.40%"

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001530"

$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "0.18%"

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0006068"

$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "0.40%"

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.006111"

$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-0.94%"

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.448"

$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-0.57%"

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.147"

$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-0.43%"

$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-0.64%"

$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "2.70%"

$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "0.84%"

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.086"

$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "5.70%"

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04252"

$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-0.14%"

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001221"

$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "0.08%"

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.003798"

$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-20.28%"

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0001200"

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-0.05%"

$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "3.31%"

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.03863"

$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "2.49%"

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1112"

$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "1.74%"

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.004083"

$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-34.29%"

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.01622"

$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "15.16%"

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.002216"

$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "0.68%"

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005155"

$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-3.05%"

$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-0.05%"

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.1355"

$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-43.89%"

$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "23.68%"

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00002100"

$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-0.05%"

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0002000"

$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-0.05%"
